$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for Harvey Antonio Valle Jimenez; the row below (Francisco
# Javier Munoz Vanegas, which carries the "last row" bottom-border style)
# shifts up into row 17 and keeps its own formatting.
$ws.Rows("17:17").Delete()

# Re-purpose that (now row 17) entry for Harvey Antonio Valle Jimenez with
# an updated "Valor Mora" amount.
$ws.Range("C17").Value = "73195757"
$ws.Range("D17").Value = "HARVEY ANTONIO VALLE JIMENEZ"
$ws.Range("E17").Value = "2508"
$ws.Range("G17").Value = 1432500

# Update the "Periodo Mora" text for the remaining worker row too.
$ws.Range("E16").Value = "2508"

# Update the summary figures.
$ws.Range("E11").Value = 114600
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 1

# Column D auto-fit shrinks now that the longest name was removed.
$ws.Columns("D").ColumnWidth = 32.3
